# Applies updated imbalance-model metrics (Kmeans(5) and Kmeans(10) segments)
# for both "test" and "valid" datasets, per the split_data test-case fixes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# Helper: write a plain-text value into a cell without Excel re-interpreting
# strings that look numeric/percent-like (e.g. "24.61...%") as a number.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 4: Kmeans(5) / test
$ws.Range("E4").Value = "[[56819    38]" + $nl + " [   24    81]]"
$ws.Range("F4").Value = 0.9989115550718023
$ws.Range("G4").Value = 0.680672268907563
$ws.Range("H4").Value = 0.7714285714285715
$ws.Range("I4").Value = 0.7232142857142857
$ws.Range("L4").Value = 3194.11
Set-TextValue $ws.Range("M4") "24.611594790288542%"
$ws.Range("N4").Value = 3

# Row 5: Kmeans(5) / valid
$ws.Range("E5").Value = "[[45464    29]" + $nl + " [   19    57]]"
$ws.Range("F5").Value = 0.9989466523294345
$ws.Range("G5").Value = 0.6627906976744186
$ws.Range("H5").Value = 0.75
$ws.Range("I5").Value = 0.7037037037037037
$ws.Range("L5").Value = 272.23
Set-TextValue $ws.Range("M5") "3.8714898046833897%"
$ws.Range("N5").Value = 3

# Row 8: Kmeans(10) / test
$ws.Range("E8").Value = "[[55163  1694]" + $nl + " [   17    88]]"
$ws.Range("F8").Value = 0.9699624310944138
$ws.Range("G8").Value = 0.04938271604938271
$ws.Range("H8").Value = 0.8380952380952381
$ws.Range("I8").Value = 0.09326974032856386
$ws.Range("L8").Value = 31273.06
Set-TextValue $ws.Range("M8") "240.96849531555927%"
$ws.Range("N8").Value = 6

# Row 9: Kmeans(10) / valid
$ws.Range("E9").Value = "[[44267  1226]" + $nl + " [   12    64]]"
$ws.Range("F9").Value = 0.9728324079966644
$ws.Range("G9").Value = 0.04961240310077519
$ws.Range("H9").Value = 0.8421052631578947
$ws.Range("I9").Value = 0.09370424597364568
$ws.Range("L9").Value = 16635.83
Set-TextValue $ws.Range("M9") "236.5846755958053%"
$ws.Range("N9").Value = 6
